# Generate Report for Handback
#
# The localization-status report is regenerated once the de-de/zh-cn
# handback packages are in sync with en-US again:
#   - Status moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" (Overview + both language sheets).
#   - The zh-cn / de-de "Latest Handback DateTime" stamps advance to the
#     new handback run time.
#   - The stale "version of handback file is not the latest" error detail
#     is cleared now that the handback is current.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# Overview sheet: Status columns for zh-cn (E2) and de-de (F2)
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# zh-cn sheet: Status, Latest Handback DateTime, Error Detail
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-17 02:52:13"
$zhcn.Range("P2").Value = ""

# de-de sheet: Status, Latest Handback DateTime, Error Detail
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-17 02:52:20"
$dede.Range("P2").Value = ""

# Re-fit the columns whose content width changed (Status got longer,
# Error Detail got shorter/empty).
$overview.Columns.Item(5).EntireColumn.AutoFit() | Out-Null
$overview.Columns.Item(6).EntireColumn.AutoFit() | Out-Null

$zhcn.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$zhcn.Columns.Item(16).EntireColumn.AutoFit() | Out-Null

$dede.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$dede.Columns.Item(16).EntireColumn.AutoFit() | Out-Null
